# Add a fourth measurement series ("WU" / "WU Error") to the wattHours
# resistive-load worksheet and wire it into the existing scatter chart,
# switching the chart's existing series over to smooth lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. New header cells (I4:J4)
# ---------------------------------------------------------------------
$ws.Range("I4").Value = "WU"
$ws.Range("J4").Value = "WU Error"

# ---------------------------------------------------------------------
# 2. New data columns I (WU readings) and J (WU Error, % vs PLM calc)
# ---------------------------------------------------------------------
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

$ws.Range("I6").Value = 9.3000000000000007
$ws.Range("J6").Formula = "=ABS(I6-D6)/D6*100"

$ws.Range("I7").Value = 18.899999999999999
$ws.Range("J7").Formula = "=ABS(I7-D7)/D7*100"

$ws.Range("I8").Value = 47.5
$ws.Range("J8").Formula = "=ABS(I8-D8)/D8*100"

$ws.Range("I9").Value = 94.7
$ws.Range("J9").Formula = "=ABS(I9-D9)/D9*100"

$ws.Range("I10").Value = 189.5
$ws.Range("J10").Formula = "=ABS(I10-D10)/D10*100"

$ws.Range("I11").Value = 285
$ws.Range("J11").Formula = "=ABS(I11-D11)/D11*100"

$ws.Range("I12").Value = 380
$ws.Range("J12").Formula = "=ABS(I12-D12)/D12*100"

$ws.Range("I13").Value = 475
$ws.Range("J13").Formula = "=ABS(I13-D13)/D13*100"

$ws.Range("I14").Value = 570
$ws.Range("J14").Formula = "=ABS(I14-D14)/D14*100"

# ---------------------------------------------------------------------
# 3. Chart updates: smooth the existing 3 series and plot the new one
# ---------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$sc = $chart.SeriesCollection()

for ($i = 1; $i -le $sc.Count; $i++) {
    $existing = $sc.Item($i)
    $existing.Smooth = $true
    $existing.MarkerStyle = -4142
}

$newSeries = $sc.NewSeries()
$newSeries.XValues = "=Sheet1!`$A`$5:`$A`$14"
$newSeries.Values = "=Sheet1!`$I`$5:`$I`$14"
$newSeries.Smooth = $true
$newSeries.MarkerStyle = -4142

# ---------------------------------------------------------------------
# 4. Restore the (incidental) selection left behind by the edit
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("U22").Select()
